$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Received Date:" value (D1) ---
$ws.Range("D1").Value = "09-04-2019 12:18:09"

# --- Row 5 (regression data set #1) ---
$ws.Range("A5").Value = "09-04-2019 12:18:09"
$ws.Range("B5").Value = "FT19040900002"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "20190408110339126"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = "FT19040900002"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "00112233440042"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "Regression AprilOne"
$ws.Range("H5").Value = "Advans Bank"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1600.00"
$ws.Range("I5").Style = "Normal"

# --- Row 6 (regression data set #2) ---
$ws.Range("A6").Value = "09-04-2019 12:18:09"
$ws.Range("B6").Value = "FT1904090000X"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "20190408110339126"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = "FT1904090000X"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "00112233440042"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "Regression AprilOne"
$ws.Range("H6").Value = "Advans Bank"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1600.00"
$ws.Range("I6").Style = "Normal"

# --- Row 7 (regression data set #3) ---
$ws.Range("A7").Value = "09-04-2019 12:18:09"
$ws.Range("B7").Value = "FT19040900002"
$ws.Range("C7").Value = "2019040811033912X"
$ws.Range("D7").Value = "FT1904090000Y"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "00112233440042"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "Regression AprilOne"
$ws.Range("H7").Value = "Advans Bank"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "1600.00"
$ws.Range("I7").Style = "Normal"

# --- Row 8 (new regression data set #4) ---
$ws.Range("A8").Value = "09-04-2019 12:18:09"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "20190408110339126"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = "FT1904090000Z"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "00112233440042"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "Regression AprilOne"
$ws.Range("H8").Value = "Advans Bank"
$ws.Range("I8").Value = "'1610"

# Row 8 previously held only row-height formatting with no cell content;
# restore the natural (unset) row height now that it carries data.
$ws.Rows.Item(8).AutoFit()

# Update the saved cursor/selection position to match the authored file.
[void]$ws.Range("B6").Select()
